# Fixed symbol conversion issues pointed out by Zainab
# cleaned up dataset
#
# The original dataset had a handful of shared strings whose special
# characters (section sign "§", curly apostrophe, en-dash) got mangled
# during an earlier import (mojibake like "Â§" and "â€“"/"â€™"). This
# script rewrites those specific cells with cleaned-up text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "Protects creators'™ rights in literary, musical, and artistic works."
$ws.Range("B18").Value = "Criminal Code, ss. 406 & 409"
$ws.Range("B4").Value  = "DMCA, 17 U.S.C. § 512"
$ws.Range("B5").Value  = "Hatch-Waxman Act, 21 U.S.C. § 355"
$ws.Range("B6").Value  = "Patent Term Restoration, 35 U.S.C. § 154"
$ws.Range("B7").Value  = "Paragraph IV Certifications, 21 U.S.C. § 355(j)(2)(A)(vii)"
$ws.Range("B8").Value  = "Design Patents, 35 U.S.C. § 171"
$ws.Range("B9").Value  = "Defend Trade Secrets Act, 18 U.S.C. § 1836"

# Restore the cursor/selection to where the author left off reviewing the fix
[void]$ws.Range("I13").Select()
